$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3503.0344
$ws.Range("J17").Value = 3521.8928
$ws.Range("L17").Value = 10565.6784
$ws.Range("N17").Value = -10901.6784
# Row 32
$ws.Range("H32").Value = 990.75
$ws.Range("I32").Value = 978.5
$ws.Range("K32").Value = 978.5
$ws.Range("M32").Value = -652.5
# Row 51
$ws.Range("H51").Value = 5199.95
# Row 69
$ws.Range("H69").Value = 9941.429
$ws.Range("J69").Value = 9941.429
$ws.Range("L69").Value = 29824.287
$ws.Range("N69").Value = -31572.287
# Row 72
$ws.Range("H72").Value = 9941.429
$ws.Range("J72").Value = 9941.429
$ws.Range("L72").Value = 89472.861
$ws.Range("N72").Value = -98208.861
# Row 98
$ws.Range("H98").Value = 1681
$ws.Range("I98").Value = 796.3333
$ws.Range("K98").Value = 796.3333
$ws.Range("M98").Value = 701.6667
# Row 113
$ws.Range("H113").Value = 4509.375
$ws.Range("J113").Value = 2038.0667
$ws.Range("L113").Value = 2038.0667
$ws.Range("N113").Value = -8546.066699999999
# Row 122
$ws.Range("H122").Value = 1681
$ws.Range("I122").Value = 796.3333
$ws.Range("K122").Value = 2388.9999
$ws.Range("M122").Value = 61.0001000000002
# Row 131
$ws.Range("H131").Value = 6359.6
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = ""
# Row 137
$ws.Range("H137").Value = 1635331.8
$ws.Range("I137").Value = 1105.1111
$ws.Range("J137").Value = 7938777.5
$ws.Range("K137").Value = 3315.3333
$ws.Range("L137").Value = 23816332.5
$ws.Range("M137").Value = -765.3333000000002
$ws.Range("N137").Value = -23821432.5
# Row 138
$ws.Range("H138").Value = 2384.34
$ws.Range("I138").Value = 896.1053000000001
$ws.Range("J138").Value = 2733.4321
$ws.Range("K138").Value = 2688.3159
$ws.Range("L138").Value = 8200.2963
$ws.Range("M138").Value = 2451.6841
$ws.Range("N138").Value = -18480.2963

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 21638296
$ws.Range("I32").Value = 24225996
$ws.Range("J32").Value = 7147167.5
$ws.Range("K32").Value = 24225996
$ws.Range("L32").Value = 7147167.5
$ws.Range("M32").Value = -24225709
$ws.Range("N32").Value = -7147741.5
# Row 45
$ws.Range("H45").Value = 3283.5
$ws.Range("I45").Value = 2197.8
$ws.Range("K45").Value = 2197.8
$ws.Range("M45").Value = -1820.8
# Row 61
$ws.Range("H61").Value = 2977.5334
$ws.Range("I61").Value = 2575.0527
$ws.Range("K61").Value = 2575.0527
$ws.Range("M61").Value = -2363.0527
# Row 74
$ws.Range("H74").Value = 2459.5
$ws.Range("I74").Value = 2231.6453
$ws.Range("J74").Value = 3244.3333
$ws.Range("K74").Value = 2231.6453
$ws.Range("L74").Value = 3244.3333
$ws.Range("M74").Value = -1357.6453
$ws.Range("N74").Value = -4992.3333
# Row 77
$ws.Range("H77").Value = 2459.5
$ws.Range("I77").Value = 2231.6453
$ws.Range("J77").Value = 3244.3333
$ws.Range("K77").Value = 11158.2265
$ws.Range("L77").Value = 16221.6665
$ws.Range("M77").Value = -6790.226500000001
$ws.Range("N77").Value = -24957.6665
# Row 122
$ws.Range("H122").Value = 3284.0908
$ws.Range("I122").Value = 2450.5386
$ws.Range("K122").Value = 7351.6158
$ws.Range("M122").Value = -4901.6158
# Row 136
$ws.Range("H136").Value = 2977.5334
$ws.Range("I136").Value = 2575.0527
$ws.Range("K136").Value = 7725.158100000001
$ws.Range("M136").Value = -5175.158100000001
# Row 139
$ws.Range("H139").Value = 80665.336
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").Value = ""

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 314
$ws.Range("I22").Value = 299.6
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 299.6
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -126.6
$ws.Range("N22").Value = -696
# Row 97
$ws.Range("H97").Value = 36885.6
$ws.Range("J97").Value = 38000
$ws.Range("L97").Value = 38000
$ws.Range("N97").Value = -39982

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 100000010
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = ""
# Row 22
$ws.Range("H22").Value = 270.75
$ws.Range("I22").Value = 249.90909
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 249.90909
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 100.09091
$ws.Range("N22").Value = -1200
# Row 31
$ws.Range("H31").Value = 3630.818
$ws.Range("I31").Value = 1229.4
$ws.Range("J31").Value = 6512.52
$ws.Range("K31").Value = 1229.4
$ws.Range("L31").Value = 6512.52
$ws.Range("M31").Value = -934.4000000000001
$ws.Range("N31").Value = -7102.52
# Row 34
$ws.Range("H34").Value = 3630.818
$ws.Range("I34").Value = 1229.4
$ws.Range("J34").Value = 6512.52
$ws.Range("K34").Value = 1229.4
$ws.Range("L34").Value = 6512.52
$ws.Range("M34").Value = -1027.4
$ws.Range("N34").Value = -6916.52
# Row 134
$ws.Range("H134").Value = 2486.111
$ws.Range("I134").Value = 2283.4
$ws.Range("K134").Value = 6850.200000000001
$ws.Range("M134").Value = -4315.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 251862260
$ws.Range("I4").Value = 214011360
$ws.Range("J4").Value = 319993900
$ws.Range("K4").Value = 642034080
$ws.Range("L4").Value = 959981700
$ws.Range("M4").Value = -642033968
$ws.Range("N4").Value = -959981924
# Row 7
$ws.Range("H7").Value = 281.75
$ws.Range("I7").Value = 148.42857
$ws.Range("J7").Value = 468.4
$ws.Range("K7").Value = 445.28571
$ws.Range("L7").Value = 1405.2
$ws.Range("M7").Value = -333.28571
$ws.Range("N7").Value = -1629.2
# Row 26
$ws.Range("H26").Value = 124.44444
$ws.Range("I26").Value = 44
$ws.Range("J26").Value = 225
$ws.Range("K26").Value = 132
$ws.Range("L26").Value = 675
$ws.Range("M26").Value = 156
$ws.Range("N26").Value = -1251
# Row 86
$ws.Range("H86").Value = 774.6667
$ws.Range("I86").Value = 399
$ws.Range("J86").Value = 962.5
$ws.Range("K86").Value = 1197
$ws.Range("L86").Value = 2887.5
$ws.Range("M86").Value = -11
$ws.Range("N86").Value = -5259.5
# Row 89
$ws.Range("H89").Value = 774.6667
$ws.Range("I89").Value = 399
$ws.Range("J89").Value = 962.5
$ws.Range("K89").Value = 3591
$ws.Range("L89").Value = 8662.5
$ws.Range("M89").Value = 2337
$ws.Range("N89").Value = -20518.5
# Row 109
$ws.Range("H109").Value = 3108.4333
$ws.Range("I109").Value = 1985.1111
$ws.Range("J109").Value = 3589.8572
$ws.Range("K109").Value = 5955.3333
$ws.Range("L109").Value = 10769.5716
$ws.Range("M109").Value = -4915.3333
$ws.Range("N109").Value = -12849.5716

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 9796
$ws.Range("I5").Value = 9796
$ws.Range("K5").Value = 9796
$ws.Range("M5").Value = -9684
# Row 70
$ws.Range("H70").Value = 29576
$ws.Range("I70").Value = 67397.60000000001
$ws.Range("J70").Value = 5937.5
$ws.Range("K70").Value = 67397.60000000001
$ws.Range("L70").Value = 5937.5
$ws.Range("M70").Value = -67127.60000000001
$ws.Range("N70").Value = -6477.5
# Row 73
$ws.Range("H73").Value = 29576
$ws.Range("I73").Value = 67397.60000000001
$ws.Range("J73").Value = 5937.5
$ws.Range("K73").Value = 67397.60000000001
$ws.Range("L73").Value = 5937.5
$ws.Range("M73").Value = -66461.60000000001
$ws.Range("N73").Value = -7809.5
# Row 97
$ws.Range("H97").Value = 724.8
$ws.Range("I97").Value = 479.5
$ws.Range("J97").Value = 1706
$ws.Range("K97").Value = 479.5
$ws.Range("L97").Value = 1706
$ws.Range("M97").Value = 16.5
$ws.Range("N97").Value = -2698
# Row 113
$ws.Range("H113").Value = 18405.5
$ws.Range("J113").Value = 26505.5
$ws.Range("L113").Value = 26505.5
$ws.Range("N113").Value = -30845.5
# Row 132
$ws.Range("H132").Value = 4349.8
$ws.Range("I132").Value = 3687.25
$ws.Range("K132").Value = 11061.75
$ws.Range("M132").Value = -8531.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 2000.4
$ws.Range("I10").Value = 2334
$ws.Range("K10").Value = 2334
$ws.Range("M10").Value = -2194
# Row 45
$ws.Range("H45").Value = 29041
$ws.Range("I45").Value = 29041
$ws.Range("K45").Value = 29041
$ws.Range("M45").Value = -28634
# Row 82
$ws.Range("H82").Value = 2427.16
$ws.Range("I82").Value = 2257.7222
$ws.Range("K82").Value = 2257.7222
$ws.Range("M82").Value = -1896.7222
# Row 85
$ws.Range("H85").Value = 2427.16
$ws.Range("I85").Value = 2257.7222
$ws.Range("K85").Value = 2257.7222
$ws.Range("M85").Value = -1009.7222
# Row 93
$ws.Range("H93").Value = 3785.7144
$ws.Range("I93").Value = 3000
$ws.Range("J93").Value = 3916.6667
$ws.Range("K93").Value = 3000
$ws.Range("L93").Value = 3916.6667
$ws.Range("M93").Value = -1752
$ws.Range("N93").Value = -6412.6667
# Row 122
$ws.Range("H122").Value = 12377
$ws.Range("I122").Value = 12964.238
$ws.Range("J122").Value = 10321.667
$ws.Range("K122").Value = 38892.714
$ws.Range("L122").Value = 30965.001
$ws.Range("M122").Value = -36442.714
$ws.Range("N122").Value = -35865.001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1490.1786
$ws.Range("I136").Value = 1184.1666
$ws.Range("J136").Value = 3326.25
$ws.Range("K136").Value = 3552.4998
$ws.Range("L136").Value = 9978.75
$ws.Range("M136").Value = -1002.4998
$ws.Range("N136").Value = -15078.75
